$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column A, shifting B:F left to A:E
$ws.Columns.Item(1).Delete()
